# Update column G ("K") values for rows 2-11 on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 8
    3  = 2
    4  = 3
    5  = 0
    6  = 6
    7  = 3
    8  = 5
    9  = 1
    10 = 2
    11 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
